$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows (33-39) appended to the "Data" sheet, matching the existing
# A1:O32 layout (col A = date/time serial, cols B:O = numeric counters).
$newRows = @(
    @(45742.222743055558, 10, 6, 280, 454, 430, 478, 3437, 478, 2026, 208, 418, 30, 3941, 5180),
    @(45745.18310185185,  10, 6, 286, 466, 441, 478, 3437, 478, 2026, 208, 418, 30, 3991, 5286),
    @(45745.185428240744, 10, 6, 286, 466, 441, 478, 3437, 478, 2026, 208, 418, 30, 3991, 5286),
    @(45749.147685185184, 10, 6, 299, 474, 449, 534, 4036, 534, 2556, 255, 494, 30, 4310, 5373),
    @(45749.251122685186, 10, 6, 299, 474, 449, 534, 4036, 534, 2556, 255, 494, 30, 4310, 5373),
    @(45749.331574074073, 10, 6, 299, 474, 449, 534, 4036, 534, 2556, 255, 494, 30, 4314, 5373),
    @(45750.121944444443, 10, 6, 299, 474, 449, 534, 4036, 534, 2570, 257, 494, 30, 4314, 5402)
)

$startRow = 33
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $values = $newRows[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $values[$c]
    }
}
